# Generate Report for Handback
#
# This models a localization "handback" run that:
#   1. Flips the Status cells from "Ready for handoff" to
#      "Handed back: in sync with en-US" (shared by the Overview sheet's
#      zh-cn/de-de status columns and both language sheets' Status column).
#   2. Fills in the "Latest Target File" / "Latest Handback File" /
#      "Latest Handback DateTime" columns on the zh-cn and de-de sheets,
#      now that a handback has actually happened for both documents.
#   3. Turns the newly-populated "Latest Target File" cells into hyperlinks
#      (same look/target as the existing "Source File Name" links).
#   4. Widens a few columns so the new, longer text isn't clipped.

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = $true
$hyperlinkColor = 15570276   # OLE (BGR) form of RGB(0x64,0x95,0xED) == FF6495ED

function Set-HandbackStatus($ws, $cellRef) {
    $ws.Range($cellRef).Value2 = "Handed back: in sync with en-US"
}

function Add-TargetFileLink($ws, $cellRef, $url, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $displayText)
    $cell = $ws.Range($cellRef)
    $cell.Font.Underline = $hyperlinkUnderline
    $cell.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet: both language Status columns move to "handed back".
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
Set-HandbackStatus $overview "E2"
Set-HandbackStatus $overview "F2"
Set-HandbackStatus $overview "E3"
Set-HandbackStatus $overview "F3"

$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
Set-HandbackStatus $zhcn "C2"
Set-HandbackStatus $zhcn "C3"

$md1Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a54c7f933a68c08f59853ff5d8e2f400181597af/e2e/3a01482a-8d86-4175-8b83-9900b91296e2.md"
$md2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a54c7f933a68c08f59853ff5d8e2f400181597af/e2e/d6200de6-3eaf-4b4a-8a06-4760996ea8ac.md"

# Row 2 (3a01482a...)
Add-TargetFileLink $zhcn "I2" $md1Url "3a01482a-8d86-4175-8b83-9900b91296e2.md"
$zhcn.Range("J2").Value2 = "3a01482a-8d86-4175-8b83-9900b91296e2.cc71d281864cdf239749586d845ce48967ad4924.zh-cn.xlf"
$zhcn.Range("K2").Value2 = "2016-08-28 02:29:54"

# Row 3 (d6200de6...)
Add-TargetFileLink $zhcn "I3" $md2Url "d6200de6-3eaf-4b4a-8a06-4760996ea8ac.md"
$zhcn.Range("J3").Value2 = "d6200de6-3eaf-4b4a-8a06-4760996ea8ac.915cb6434337f34921815e02a316133e8e0b24d0.zh-cn.xlf"
$zhcn.Range("K3").Value2 = "2016-08-28 02:29:54"

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
Set-HandbackStatus $dede "C2"
Set-HandbackStatus $dede "C3"

# Row 2 (3a01482a...)
Add-TargetFileLink $dede "I2" $md1Url "3a01482a-8d86-4175-8b83-9900b91296e2.md"
$dede.Range("J2").Value2 = "3a01482a-8d86-4175-8b83-9900b91296e2.cc71d281864cdf239749586d845ce48967ad4924.de-de.xlf"
$dede.Range("K2").Value2 = "2016-08-28 02:30:07"

# Row 3 (d6200de6...)
Add-TargetFileLink $dede "I3" $md2Url "d6200de6-3eaf-4b4a-8a06-4760996ea8ac.md"
$dede.Range("J3").Value2 = "d6200de6-3eaf-4b4a-8a06-4760996ea8ac.915cb6434337f34921815e02a316133e8e0b24d0.de-de.xlf"
$dede.Range("K3").Value2 = "2016-08-28 02:30:07"

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15

Write-Output "Handback report generated."
